$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.446.94"
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.712.43"
$ws.Range("E3").Value = "  -1.65%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.08"
$ws.Range("E5").Value = "  -2.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.83"
$ws.Range("E6").Value = "  -1.39%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.590"
$ws.Range("E8").Value = "  -2.31%  "
$ws.Range("E9").Value = "  -2.59%  "
$ws.Range("E10").Value = "  +1.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.57"
$ws.Range("E11").Value = "  -1.59%  "
$ws.Range("E12").Value = "  -3.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.196.31"
$ws.Range("E13").Value = "  -1.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.43"
$ws.Range("E14").Value = "  -1.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.360.16"
$ws.Range("E15").Value = "  -0.52%  "
$ws.Range("E16").Value = "  -2.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.719.62"
$ws.Range("E17").Value = "  -1.52%  "
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.66"
$ws.Range("E19").Value = "  -4.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "351.28"
$ws.Range("E20").Value = "  -1.36%  "
$ws.Range("E21").Value = "  -4.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("E23").Value = "  -4.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.36"
$ws.Range("E24").Value = "  -1.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.169"
$ws.Range("E25").Value = "  -1.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  -4.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0888"
$ws.Range("E28").Value = "  -2.20%  "
$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.36"
$ws.Range("E29").Value = "  +9.30%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.95"
$ws.Range("E30").Value = "  -0.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.17"
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.84"
$ws.Range("E32").Value = "  -2.06%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  -2.28%  "
$ws.Range("E36").Value = "  -1.85%  "
$ws.Range("E37").Value = "  -2.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "345.59"
$ws.Range("E38").Value = "  +0.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.964"
$ws.Range("E39").Value = "  -4.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.12"
$ws.Range("E40").Value = "  -3.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.04"
$ws.Range("E41").Value = "  -3.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.45"
$ws.Range("E42").Value = "  -1.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.40"
$ws.Range("E43").Value = "  -1.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.73"
$ws.Range("E44").Value = "  -3.19%  "
$ws.Range("E45").Value = "  -3.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.622"
$ws.Range("E46").Value = "  -1.57%  "
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "131.68"
$ws.Range("E48").Value = "  -3.18%  "
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0983"
$ws.Range("E51").Value = "  -3.85%  "
